# Commit: "Fruta / hortaliza, semanal" — weekly refresh of the daily
# Cilantro logic-table: a new record (week) is inserted at the top of the
# data block (row 76) and every subsequent record shifts down by one row,
# extending the table from A1:R150 to A1:R151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 76..150 down to 77..151, opening up a blank row 76.
$ws.Rows.Item(76).Insert()

# Populate the newly opened row with the new weekly record.
$ws.Cells.Item(76, 1).Value  = 8
$ws.Cells.Item(76, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76, 3).Value  = "Coquimbo"
$ws.Cells.Item(76, 4).Value  = 44740
$ws.Cells.Item(76, 5).Value  = 4
$ws.Cells.Item(76, 6).Value  = 100112040
$ws.Cells.Item(76, 7).Value  = "Cilantro"
$ws.Cells.Item(76, 8).Value  = "Sin especificar"
$ws.Cells.Item(76, 9).Value  = "Primera"
$ws.Cells.Item(76, 10).Value = 2800
$ws.Cells.Item(76, 11).Value = 1300
$ws.Cells.Item(76, 12).Value = 1500
$ws.Cells.Item(76, 13).Value = 1400
$ws.Cells.Item(76, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(76, 16).Value = 933
$ws.Cells.Item(76, 17).Value = 1.5
$ws.Cells.Item(76, 18).Value = "Hortaliza"
